$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy row 127 formatting down to new rows 128 and 129 (adds styled, empty rows) ---
$ws.Range("A127:V127").Copy()
$ws.Range("A128:V128").PasteSpecial(-4122)
$ws.Range("A129:V129").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Swap/update existing rows per diff ---
# Row 11
$ws.Range("F11").Value = 'Opava'
$ws.Range("H11").Value = 'Varnsdorf'
$ws.Range("J11").Value = 1.93
$ws.Range("L11").Value = 1.87
$ws.Range("M11").Value = '28/07/2023 17:51'
$ws.Range("N11").Value = 3.51
$ws.Range("P11").Value = 3.67
$ws.Range("Q11").Value = '28/07/2023 17:51'
$ws.Range("R11").Value = 3.46
$ws.Range("T11").Value = 4
$ws.Range("U11").Value = '28/07/2023 17:51'
$ws.Range("V11").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/opava-varnsdorf/ELFSOQvo/'

# Row 12
$ws.Range("F12").Value = 'Taborsko'
$ws.Range("H12").Value = 'Sparta Prague B'
$ws.Range("J12").Value = 1.79
$ws.Range("L12").Value = 1.91
$ws.Range("M12").Value = '28/07/2023 17:50'
$ws.Range("N12").Value = 3.62
$ws.Range("P12").Value = 3.79
$ws.Range("Q12").Value = '28/07/2023 17:56'
$ws.Range("R12").Value = 3.89
$ws.Range("T12").Value = 3.71
$ws.Range("U12").Value = '28/07/2023 17:56'
$ws.Range("V12").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/taborsko-sparta-prague/ri4VrnWG/'

# Row 84
$ws.Range("F84").Value = 'Kromeriz'
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 'Vlasim'
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 2.67
$ws.Range("L84").Value = 3.1
$ws.Range("M84").Value = '30/09/2023 15:47'
$ws.Range("N84").Value = 3.28
$ws.Range("P84").Value = 3.57
$ws.Range("Q84").Value = '30/09/2023 15:47'
$ws.Range("R84").Value = 2.37
$ws.Range("T84").Value = 2.2
$ws.Range("U84").Value = '30/09/2023 15:47'
$ws.Range("V84").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/kromeriz-vlasim/lWoUwF6A/'

# Row 85
$ws.Range("F85").Value = 'Lisen'
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 'Taborsko'
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1.95
$ws.Range("L85").Value = 2.16
$ws.Range("M85").Value = '30/09/2023 15:48'
$ws.Range("N85").Value = 3.3
$ws.Range("P85").Value = 3.31
$ws.Range("Q85").Value = '30/09/2023 15:58'
$ws.Range("R85").Value = 3.48
$ws.Range("T85").Value = 3.42
$ws.Range("U85").Value = '30/09/2023 15:48'
$ws.Range("V85").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/lisen-taborsko/KdsQvei4/'

# Row 115
$ws.Range("F115").Value = 'Lisen'
$ws.Range("H115").Value = 'Vyskov'
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 2.54
$ws.Range("L115").Value = 2.75
$ws.Range("M115").Value = '04/11/2023 13:56'
$ws.Range("N115").Value = 3.22
$ws.Range("P115").Value = 2.99
$ws.Range("Q115").Value = '04/11/2023 13:56'
$ws.Range("R115").Value = 2.58
$ws.Range("T115").Value = 2.79
$ws.Range("U115").Value = '04/11/2023 13:56'
$ws.Range("V115").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/lisen-mfk-vyskov/n5xcbtlf/'

# Row 116
$ws.Range("F116").Value = 'Vlasim'
$ws.Range("H116").Value = 'Sparta Prague B'
$ws.Range("I116").Value = 1
$ws.Range("J116").Value = 1.61
$ws.Range("L116").Value = 1.62
$ws.Range("M116").Value = '04/11/2023 13:52'
$ws.Range("N116").Value = 4
$ws.Range("P116").Value = 4.37
$ws.Range("Q116").Value = '04/11/2023 13:52'
$ws.Range("R116").Value = 4.33
$ws.Range("T116").Value = 4.75
$ws.Range("U116").Value = '04/11/2023 13:57'
$ws.Range("V116").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/vlasim-sparta-prague/GpUDfbmJ/'

# Row 117
$ws.Range("F117").Value = 'Kromeriz'
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 'Brno'
$ws.Range("I117").Value = 3
$ws.Range("J117").Value = 3.41
$ws.Range("L117").Value = 4.19
$ws.Range("M117").Value = '04/11/2023 13:52'
$ws.Range("N117").Value = 3.43
$ws.Range("P117").Value = 3.85
$ws.Range("Q117").Value = '04/11/2023 13:52'
$ws.Range("R117").Value = 1.93
$ws.Range("T117").Value = 1.79
$ws.Range("U117").Value = '04/11/2023 13:52'
$ws.Range("V117").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/kromeriz-brno/tSy1c030/'

# Row 118
$ws.Range("F118").Value = 'Opava'
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 'Prostejov'
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = 1.49
$ws.Range("L118").Value = 1.55
$ws.Range("M118").Value = '04/11/2023 13:56'
$ws.Range("N118").Value = 4.15
$ws.Range("P118").Value = 4.09
$ws.Range("Q118").Value = '04/11/2023 13:56'
$ws.Range("R118").Value = 5.25
$ws.Range("T118").Value = 6.04
$ws.Range("V118").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/opava-prostejov/8EwgaMYm/'

# Row 124
$ws.Range("F124").Value = 'Vyskov'
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 'Dukla Prague'
$ws.Range("J124").Value = 2.33
$ws.Range("L124").Value = 2.01
$ws.Range("M124").Value = '11/11/2023 09:51'
$ws.Range("N124").Value = 3.6
$ws.Range("P124").Value = 3.59
$ws.Range("Q124").Value = '11/11/2023 09:51'
$ws.Range("R124").Value = 2.59
$ws.Range("T124").Value = 3.54
$ws.Range("U124").Value = '11/11/2023 09:51'
$ws.Range("V124").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/mfk-vyskov-dukla-prague/KWFwBJXa/'

# Row 125
$ws.Range("F125").Value = 'Chrudim'
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 'Vlasim'
$ws.Range("J125").Value = 2.51
$ws.Range("L125").Value = 2.33
$ws.Range("M125").Value = '11/11/2023 10:14'
$ws.Range("N125").Value = 3.3
$ws.Range("P125").Value = 3.46
$ws.Range("Q125").Value = '11/11/2023 10:05'
$ws.Range("R125").Value = 2.51
$ws.Range("T125").Value = 2.94
$ws.Range("U125").Value = '11/11/2023 10:14'
$ws.Range("V125").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/chrudim-vlasim/x8ArAwm6/'

# --- Append two new match rows (128, 129) ---
# Row 128
$ws.Range("A128").Value = 127
$ws.Range("B128").Value = 'czech-republic'
$ws.Range("C128").Value = 'fnl'
$ws.Range("D128").Value = '2023-2024'
$ws.Range("E128").Value = 45242.4375
$ws.Range("F128").Value = 'Sparta Prague B'
$ws.Range("G128").Value = 2
$ws.Range("H128").Value = 'Taborsko'
$ws.Range("I128").Value = 1
$ws.Range("J128").Value = 2.66
$ws.Range("K128").Value = '09/11/2023 09:13'
$ws.Range("L128").Value = 2.72
$ws.Range("M128").Value = '12/11/2023 09:50'
$ws.Range("N128").Value = 3.22
$ws.Range("O128").Value = '09/11/2023 09:13'
$ws.Range("P128").Value = 3.37
$ws.Range("Q128").Value = '12/11/2023 09:50'
$ws.Range("R128").Value = 2.47
$ws.Range("S128").Value = '09/11/2023 09:13'
$ws.Range("T128").Value = 2.54
$ws.Range("U128").Value = '12/11/2023 09:50'
$ws.Range("V128").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/sparta-prague-taborsko/Gb9n9c2C/'

# Row 129
$ws.Range("A129").Value = 128
$ws.Range("B129").Value = 'czech-republic'
$ws.Range("C129").Value = 'fnl'
$ws.Range("D129").Value = '2023-2024'
$ws.Range("E129").Value = 45242.58333333334
$ws.Range("F129").Value = 'Varnsdorf'
$ws.Range("G129").Value = 2
$ws.Range("H129").Value = 'Opava'
$ws.Range("I129").Value = 1
$ws.Range("J129").Value = 2.65
$ws.Range("K129").Value = '09/11/2023 09:13'
$ws.Range("L129").Value = 2.6
$ws.Range("M129").Value = '12/11/2023 13:42'
$ws.Range("N129").Value = 3.43
$ws.Range("O129").Value = '09/11/2023 09:13'
$ws.Range("P129").Value = 3.62
$ws.Range("Q129").Value = '12/11/2023 12:51'
$ws.Range("R129").Value = 2.36
$ws.Range("S129").Value = '09/11/2023 09:13'
$ws.Range("T129").Value = 2.53
$ws.Range("U129").Value = '12/11/2023 13:57'
$ws.Range("V129").Value = 'https://www.betexplorer.com/football/czech-republic/fnl/varnsdorf-opava/riMaI1uQ/'
